$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.340.75"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "1.832.25"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.93%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4746"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3690"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07456"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8849"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "1.869.51"
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07318"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.585"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008792"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "27.568.97"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.296"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "2.095.65"
$ws.Range("E24").Value = "  +2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.891"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.137"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.245"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08994"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7542"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.544"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.010"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05345"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01953"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.978"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.300"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.399"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5315"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.476"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4914"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.672"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06297"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
